$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("python")

# --- Row 24: add "completed" marker in column D (existing row, A/B unchanged) ---
$ws.Range("D24").Value = "completed"

# --- Row 25: add topic in column B, and "completed" marker in column D (A unchanged) ---
$ws.Range("B25").Value = "stack ,topn,sort,window recepe"
$ws.Range("D25").Value = "completed"

# --- Row 26 (new) ---
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = 45686
$ws.Range("B26").Value = "how to upload local files to dataiku folder"
$ws.Range("D26").Value = "completed"

# --- Row 27 (new, no date) ---
$ws.Range("B27").Value = "how to move folder one project to another project folder"
$ws.Range("D27").Value = "completed"

# --- Row 28 (new, no date) ---
$ws.Range("B28").Value = "how to run python recepe individual"
$ws.Range("D28").Value = "completed"

# --- Row 29 (new) ---
$ws.Range("A25").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = 45687
$ws.Range("B29").Value = "python code practice"
$ws.Range("D29").Value = "completed"

# --- Row 30 (new) ---
$ws.Range("A25").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = 45688

# --- Row 31 (new): write Saturday BEFORE "dataiku options..." text so the
# shared-string table gets the same allocation order as the source workbook ---
$ws.Range("A25").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value = 45689
$ws.Range("B31").Value = "Saturday"

# Now write row 30's text (this creates the next shared string AFTER "Saturday")
$ws.Range("B30").Value = "dataiku options exploration,git hub team discussion,  "
$ws.Range("D30").Value = "completed"

# --- Row 32 (new): Sunday, with new date style (numFmtId 16 "d-mmm") ---
$ws.Range("A32").Value = 45690
$ws.Range("A32").NumberFormat = "d-mmm"
$ws.Range("B32").Value = "sunday"

# --- Row 33 (new) ---
$ws.Range("A25").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = 45691
$ws.Range("B33").Value = "Monday-client holiday"

# --- Row 34 (new) ---
$ws.Range("A25").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A34").Value = 45692
$ws.Range("B34").Value = "sql class ,variable creation,list,tuple,dict,slicing sequences,"
$ws.Range("D34").Value = "completed"

# --- Row 35 (new, no date) ---
$ws.Range("B35").Value = "creating multiple dictionaries"
$ws.Range("D35").Value = "completed"

# --- Row 36 (new) ---
$ws.Range("A25").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value = 45693
$ws.Range("B36").Value = "assignment 4(b), 5(a)"
$ws.Range("D36").Value = "completed"

# --- Row 37 (new) ---
$ws.Range("A25").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = 45694
$ws.Range("B37").Value = "dictionary creation,dictionary pandas dataframe ga create."

# --- Column B width ---
$ws.Columns("B").ColumnWidth = 54.7

# --- Selection / view ---
$ws.Range("B37").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
